# Abhinav | Add test case e3
#
# Adds two new worksheets ("e3" and "e3_shifted") after "e2_shifted",
# mirroring the layout of "e2"/"e2_shifted" but with an extra "Main subject"
# header row above the Heading row. Also restores/updates the view state
# (selections) on a handful of sheets the way the human author left them
# after finishing the edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "e3" sheet (same numbers as e2 / e2_shifted, plus a
#    "Main subject" label row above the Heading row).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$e3 = $wb.Worksheets.Add($null, $lastSheet)
$e3.Name = "e3"

$e3.Range("B3").Value = "Main subject"
$e3.Range("B4").Value = "Heading 1"
$e3.Range("C4").Value = "Heading 2"
$e3.Range("D4").Value = "Heading 3"

$data = @(
  @(12, 2, 1),
  @(23, 4, 2),
  @(34, 6, 3),
  @(45, 8, 4),
  @(56, 10, 5),
  @(67, 10, 6)
)
$row = 5
foreach ($r in $data) {
  $e3.Range("B$row").Value = $r[0]
  $e3.Range("C$row").Value = $r[1]
  $e3.Range("D$row").Value = $r[2]
  $row++
}

# ---------------------------------------------------------------------
# 2. Build the new "e3_shifted" sheet -- same content as "e3" shifted one
#    column/row over, mirroring how "e2_shifted" relates to "e2".
# ---------------------------------------------------------------------
$e3s = $wb.Worksheets.Add($null, $e3)
$e3s.Name = "e3_shifted"

$e3s.Range("C4").Value = "Main subject"
$e3s.Range("C5").Value = "Heading 1"
$e3s.Range("D5").Value = "Heading 2"
$e3s.Range("E5").Value = "Heading 3"

$row = 6
foreach ($r in $data) {
  $e3s.Range("C$row").Value = $r[0]
  $e3s.Range("D$row").Value = $r[1]
  $e3s.Range("E$row").Value = $r[2]
  $row++
}

# ---------------------------------------------------------------------
# 3. Restore per-sheet selection / scroll state to match where the author
#    left each sheet.
# ---------------------------------------------------------------------

# india_wheat: scrolled right a bit, landed on O9
$indiaWheat = $wb.Worksheets.Item("india_wheat")
$indiaWheat.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$indiaWheat.Range("O9").Select() | Out-Null

# e1: selection moved to M33
$e1 = $wb.Worksheets.Item("e1")
$e1.Range("M33").Select() | Out-Null

# e2: selection narrowed to B3:B8 (no longer the active tab)
$e2 = $wb.Worksheets.Item("e2")
$e2.Range("B3:B8").Select() | Out-Null

# e2_shifted: selection now covers B4:D10
$e2Shifted = $wb.Worksheets.Item("e2_shifted")
$e2Shifted.Range("B4:D10").Select() | Out-Null

# e3_shifted: selection parked at I17
$e3s.Range("I17").Select() | Out-Null

# e3: active tab / selection on B3 -- select this last so it ends up the
# workbook's active sheet (matches activeTab pointing at "e3").
$e3.Range("B3").Select() | Out-Null

Write-Host "Added e3 and e3_shifted; updated view state."
